{"js": "// Lattice-multiplication worksheet refresh: every cell in the single 5x3\n// table keeps its exact layout (problem line, factor line, \"----\" rule,\n// two partial-product placeholder lines separated by manual line breaks)\n// but gets a new multiplication problem. Replace each cell's text content,\n// in document order, with the new values -- the run-level formatting\n// (font size 32) and the w:br-separated line structure are preserved\n// automatically because we replace through the cell body Range, which\n// Word round-trips \"\\u000b\" back into <w:br/>.\nconst newCellText = [\n  \"29 x 65\\u000b  6    5\\u000b  ----\\u000b2|    |\\u000b9|    |\",\n  \"35 x 45\\u000b  4    5\\u000b  ----\\u000b3|    |\\u000b5|    |\",\n  \"40 x 33\\u000b  3    3\\u000b  ----\\u000b4|    |\\u000b0|    |\",\n  \"30 x 27\\u000b  2    7\\u000b  ----\\u000b3|    |\\u000b0|    |\",\n  \"72 x 65\\u000b  6    5\\u000b  ----\\u000b7|    |\\u000b2|    |\",\n  \"89 x 95\\u000b  9    5\\u000b  ----\\u000b8|    |\\u000b9|    |\",\n  \"16 x 13\\u000b  1    3\\u000b  ----\\u000b1|    |\\u000b6|    |\",\n  \"28 x 81\\u000b  8    1\\u000b  ----\\u000b2|    |\\u000b8|    |\",\n  \"92 x 25\\u000b  2    5\\u000b  ----\\u000b9|    |\\u000b2|    |\",\n  \"50 x 77\\u000b  7    7\\u000b  ----\\u000b5|    |\\u000b0|    |\",\n  \"92 x 92\\u000b  9    2\\u000b  ----\\u000b9|    |\\u000b2|    |\",\n  \"33 x 13\\u000b  1    3\\u000b  ----\\u000b3|    |\\u000b3|    |\",\n  \"87 x 97\\u000b  9    7\\u000b  ----\\u000b8|    |\\u000b7|    |\",\n  \"16 x 75\\u000b  7    5\\u000b  ----\\u000b1|    |\\u000b6|    |\",\n  \"37 x 86\\u000b  8    6\\u000b  ----\\u000b3|    |\\u000b7|    |\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nlet i = 0;\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < 3; c++) {\n    const cell = table.getCell(r, c);\n    const range = cell.body.getRange();\n    range.insertText(newCellText[i], \"Replace\");\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Lattice-multiplication worksheet refresh: every cell in the single 5x3\n# table keeps its exact layout (problem line, factor line, \"----\" rule,\n# two partial-product placeholder lines separated by manual line breaks)\n# but gets a new multiplication problem. Replace each cell's Range.Text,\n# in document order -- the run-level formatting (font size 32) is kept by\n# Word because we overwrite the existing range in place, and \"`v\" (vertical\n# tab, chr 11) round-trips back into <w:br/> just like the original manual\n# line breaks.\n$d = $word.ActiveDocument\n\n$newCellText = @(\n  \"29 x 65`v  6    5`v  ----`v2|    |`v9|    |\",\n  \"35 x 45`v  4    5`v  ----`v3|    |`v5|    |\",\n  \"40 x 33`v  3    3`v  ----`v4|    |`v0|    |\",\n  \"30 x 27`v  2    7`v  ----`v3|    |`v0|    |\",\n  \"72 x 65`v  6    5`v  ----`v7|    |`v2|    |\",\n  \"89 x 95`v  9    5`v  ----`v8|    |`v9|    |\",\n  \"16 x 13`v  1    3`v  ----`v1|    |`v6|    |\",\n  \"28 x 81`v  8    1`v  ----`v2|    |`v8|    |\",\n  \"92 x 25`v  2    5`v  ----`v9|    |`v2|    |\",\n  \"50 x 77`v  7    7`v  ----`v5|    |`v0|    |\",\n  \"92 x 92`v  9    2`v  ----`v9|    |`v2|    |\",\n  \"33 x 13`v  1    3`v  ----`v3|    |`v3|    |\",\n  \"87 x 97`v  9    7`v  ----`v8|    |`v7|    |\",\n  \"16 x 75`v  7    5`v  ----`v1|    |`v6|    |\",\n  \"37 x 86`v  8    6`v  ----`v3|    |`v7|    |\"\n)\n\n$tbl = $d.Tables.Item(1)\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $cell = $tbl.Cell($r, $c)\n    $cell.Range.Text = $newCellText[$i]\n    $i++\n  }\n}\n"}
